$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.848.63"
$ws.Range("E2").Value = "  -1.71%  "
$ws.Range("D3").Value = "3.271.62"
$ws.Range("E3").Value = "  -1.12%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.33%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.599"
$ws.Range("D8").Style = "Normal"
$ws.Range("E9").Value = "  -4.83%  "
$ws.Range("E10").Value = "  -1.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.408"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.70%  "
$ws.Range("D12").Value = "3.829.52"
$ws.Range("E12").Value = "  -1.34%  "
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.50"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.81%  "
$ws.Range("D15").Value = "67.871.28"
$ws.Range("E15").Value = "  -1.69%  "
$ws.Range("E16").Value = "  -3.20%  "
$ws.Range("D17").Value = "3.263.94"
$ws.Range("E17").Value = "  -1.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.72"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "398.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.67%  "
$ws.Range("E24").Value = "  -2.12%  "
$ws.Range("E25").Value = "  -4.70%  "
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("E27").Value = "  -2.80%  "
$ws.Range("E28").Value = "  +0.76%  "
$ws.Range("E29").Value = "  -2.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.48"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.92"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.25"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "162.59"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.41%  "
$ws.Range("E36").Value = "  -5.88%  "
$ws.Range("E37").Value = "  -1.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.808"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.50"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.31"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.63%  "
$ws.Range("D42").Value = "2.669.88"
$ws.Range("E42").Value = "  +1.17%  "
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0681"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.93%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.50%  "
$ws.Range("E45").Value = "  -8.26%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "24.59"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.60%  "
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "334.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.55%  "
$ws.Range("E48").Value = "  -4.00%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("E50").Value = "  -1.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.970"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.27%  "
